$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '64.430.59'
$ws.Range("E2").Value = '  -2.74%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '3.167.66'
$ws.Range("E3").Value = '  -4.52%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  +0.07%  '

# Row 5 - BNB
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '571.67'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -2.42%  '

# Row 6 - Solana
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '168.85'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -6.75%  '

# Row 7 - XRP
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.604'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -7.14%  '

# Row 8 - USDC
$ws.Range("E8").Value = '  -0.11%  '

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = '3.185.13'
$ws.Range("E9").Value = '  -3.98%  '

# Row 10 - Dogecoin
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.119'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -5.38%  '

# Row 11 - Toncoin
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '6.82'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.18%  '

# Row 12 - Cardano
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.389'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -2.99%  '

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '3.715.25'
$ws.Range("E13").Value = '  -4.72%  '

# Row 14 - TRON
$ws.Range("E14").Value = '  -1.57%  '

# Row 15 - WrappedBTC
$ws.Range("D15").Value = '64.507.81'
$ws.Range("E15").Value = '  -2.68%  '

# Row 16 - Avalanche
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '25.38'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -4.58%  '

# Row 17 - ShibaInu->WrappedEther
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.189.49'
$ws.Range("E17").Value = '  -4.04%  '

# Row 18 - WrappedEther->ShibaInu
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.0000157'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -3.96%  '

# Row 19 - BitcoinCash
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '417.61'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.54%  '

# Row 20 - Chainlink
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '12.88'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.81%  '

# Row 21 - Polkadot
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '5.31'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -3.59%  '

# Row 22 - Uniswap
$ws.Range("E22").Value = '  -2.95%  '

# Row 23 - Dai
$ws.Range("E23").Value = '  +0.22%  '

# Row 24 - LEO
$ws.Range("E24").Value = '  +0.09%  '

# Row 25 - Litecoin
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '69.73'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.84%  '

# Row 26 - Kaspa
$ws.Range("E26").Value = '  -0.61%  '

# Row 27 - Polygon
$ws.Range("E27").Value = '  -2.79%  '

# Row 28 - PEPE
$ws.Range("E28").Value = '  -9.52%  '

# Row 29 - InternetComputer(DFINITY)
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '8.76'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -2.99%  '

# Row 30 - Binance-PegBSC-USD
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"

# Row 31 - PancakeSwap
$ws.Range("E31").Value = '  -4.92%  '

# Row 32 - EthereumClassic
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '21.75'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.68%  '

# Row 33 - USDe
$ws.Range("E33").Value = '  -0.09%  '

# Row 34 - NEARProtocol
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '5.07'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -2.07%  '

# Row 35 - Aptos
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '6.36'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -3.96%  '

# Row 36 - Fetch.AI
$ws.Range("E36").Value = '  -5.04%  '

# Row 37 - Monero
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '155.75'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -2.72%  '

# Row 38 - ImmutableX
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.36'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -4.98%  '

# Row 39 - Stacks->Maker
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '2.705.21'
$ws.Range("E39").Value = '  -5.48%  '

# Row 40 - Maker->Stacks
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.70'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -5.36%  '

# Row 41 - Filecoin
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '4.22'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -2.55%  '

# Row 42 - EnergySwap
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '24.09'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -8.60%  '

# Row 43 - OKB
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '39.06'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.76%  '

# Row 44 - Mantle
$ws.Range("E44").Value = '  -5.24%  '

# Row 45 - Hedera
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0620'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -6.10%  '

# Row 46 - RenderToken
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '5.45'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -8.60%  '

# Row 47 - VeChain
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.0263'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -3.21%  '

# Row 48 - InjectiveProtocol
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '21.48'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -6.82%  '

# Row 49 - Bittensor
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '290.38'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -6.83%  '

# Row 50 - FirstDigitalUSD
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.10%  '

# Row 51 - Stellar
$ws.Range("E51").Value = '  -5.11%  '
